$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 13.42939298204104
$ws.Range("D2").Value = 4.869242924293983
$ws.Range("E2").Value = 13.79431582686956
$ws.Range("F2").Value = 24.53724826728914
$ws.Range("G2").Value = 29.3015829382752
$ws.Range("H2").Value = 14.33777293421344
$ws.Range("K2").Value = 14.03485835500682
$ws.Range("L2").Value = 9.199061398802536
$ws.Range("N2").Value = 16.81734443349527
$ws.Range("O2").Value = 21.91599288424075
$ws.Range("C3").Value = 13.32876118497864
$ws.Range("D3").Value = 4.827578687229076
$ws.Range("E3").Value = 13.73191295437658
$ws.Range("F3").Value = 24.53765085262209
$ws.Range("G3").Value = 29.28443232432006
$ws.Range("H3").Value = 14.38238421520306
$ws.Range("K3").Value = 13.54223955579664
$ws.Range("L3").Value = 9.200557490464773
$ws.Range("N3").Value = 16.84032109496912
$ws.Range("O3").Value = 21.97179184343366
$ws.Range("C4").Value = 13.26997333803254
$ws.Range("D4").Value = 4.801474167779489
$ws.Range("E4").Value = 13.69666690837401
$ws.Range("F4").Value = 24.54537314795261
$ws.Range("G4").Value = 29.28495756246945
$ws.Range("H4").Value = 14.41251123619463
$ws.Range("K4").Value = 13.23140681746614
$ws.Range("L4").Value = 9.20298955061503
$ws.Range("N4").Value = 16.85618402148108
$ws.Range("O4").Value = 22.01169604609113
$ws.Range("C5").Value = 13.2467944022875
$ws.Range("D5").Value = 4.790709549948375
$ws.Range("E5").Value = 13.68308710470094
$ws.Range("F5").Value = 24.55039772802352
$ws.Range("G5").Value = 29.2879473467448
$ws.Range("H5").Value = 14.42547505670809
$ws.Range("K5").Value = 13.10281935828769
$ws.Range("L5").Value = 9.204361732606557
$ws.Range("N5").Value = 16.86309045970093
$ws.Range("O5").Value = 22.02937174628736
$ws.Range("C6").Value = 13.24299316858986
$ws.Range("D6").Value = 4.788914568704493
$ws.Range("E6").Value = 13.68087980610498
$ws.Range("F6").Value = 24.55134537690534
$ws.Range("G6").Value = 29.28861126248492
$ws.Range("H6").Value = 14.42766913976507
$ws.Range("K6").Value = 13.08135718477585
$ws.Range("L6").Value = 9.204612612467432
$ws.Range("N6").Value = 16.86426399801427
$ws.Range("O6").Value = 22.03239206684321
$ws.Range("C7").Value = 13.26965756092644
$ws.Range("D7").Value = 4.801329499777003
$ws.Range("E7").Value = 13.69648058082298
$ws.Range("F7").Value = 24.54543331225405
$ws.Range("G7").Value = 29.28498665279326
$ws.Range("H7").Value = 14.41268329135413
$ws.Range("K7").Value = 13.2296801645521
$ws.Range("L7").Value = 9.203006512695742
$ws.Range("N7").Value = 16.85627537273894
$ws.Range("O7").Value = 22.01192870630623
$ws.Range("C8").Value = 13.39408667470364
$ws.Range("D8").Value = 4.85498826765129
$ws.Range("E8").Value = 13.77216937532498
$ws.Range("F8").Value = 24.53583484459908
$ws.Range("G8").Value = 29.29337306814173
$ws.Range("H8").Value = 14.35258652524987
$ws.Range("K8").Value = 13.86684300839468
$ws.Range("L8").Value = 9.199263488523806
$ws.Range("N8").Value = 16.82490299408932
$ws.Range("O8").Value = 21.93405836613293
$ws.Range("C9").Value = 13.66081870223104
$ws.Range("D9").Value = 4.955870367751936
$ws.Range("E9").Value = 13.94440192862844
$ws.Range("F9").Value = 24.576350147614
$ws.Range("G9").Value = 29.39759699855773
$ws.Range("H9").Value = 14.25648810082648
$ws.Range("K9").Value = 15.04284031279873
$ws.Range("L9").Value = 9.203904219825311
$ws.Range("N9").Value = 16.77727284170416
$ws.Range("O9").Value = 21.82632773420527
$ws.Range("C10").Value = 13.86911118404767
$ws.Range("D10").Value = 5.027084218583426
$ws.Range("E10").Value = 14.08465314729752
$ws.Range("F10").Value = 24.64223461692846
$ws.Range("G10").Value = 29.52755797630016
$ws.Range("H10").Value = 14.19920618337973
$ws.Range("K10").Value = 15.85370176986224
$ws.Range("L10").Value = 9.214573435118984
$ws.Range("N10").Value = 16.75070022795672
$ws.Range("O10").Value = 21.77484007601483
$ws.Range("C11").Value = 13.96619915753759
$ws.Range("D11").Value = 5.058797674224362
$ws.Range("E11").Value = 14.15124471002431
$ws.Range("F11").Value = 24.68000607820283
$ws.Range("G11").Value = 29.59818131756521
$ws.Range("H11").Value = 14.17605228447748
$ws.Range("K11").Value = 16.20959567604502
$ws.Range("L11").Value = 9.22099116259975
$ws.Range("N11").Value = 16.74043007344876
$ws.Range("O11").Value = 21.75746795787752
$ws.Range("C12").Value = 14.00326950441202
$ws.Range("D12").Value = 5.070704270841561
$ws.Range("E12").Value = 14.17684525637572
$ws.Range("F12").Value = 24.69542454081416
$ws.Range("G12").Value = 29.62656596936758
$ws.Range("H12").Value = 14.16770306028536
$ws.Range("K12").Value = 16.34239666229758
$ws.Range("L12").Value = 9.223645001248306
$ws.Range("N12").Value = 16.73680155815892
$ws.Range("O12").Value = 21.75176246532657
$ws.Range("C13").Value = 13.99527265640647
$ws.Range("D13").Value = 5.068144606238775
$ws.Range("E13").Value = 14.17131493879211
$ws.Range("F13").Value = 24.69205441806035
$ws.Range("G13").Value = 29.62038007810107
$ws.Range("H13").Value = 14.16948257876049
$ws.Range("K13").Value = 16.31388463925278
$ws.Range("L13").Value = 9.223063528792869
$ws.Range("N13").Value = 16.73757144846557
$ws.Range("O13").Value = 21.75295237785264
$ws.Range("C14").Value = 13.96924301603512
$ws.Range("D14").Value = 5.059779318872916
$ws.Range("E14").Value = 14.15334329391872
$ws.Range("F14").Value = 24.68125225598735
$ws.Range("G14").Value = 29.60048373516768
$ws.Range("H14").Value = 14.17535699234423
$ws.Range("K14").Value = 16.22056128550085
$ws.Range("L14").Value = 9.221205022825787
$ws.Range("N14").Value = 16.74012633596937
$ws.Range("O14").Value = 21.75698105286293
$ws.Range("C15").Value = 13.95333794834517
$ws.Range("D15").Value = 5.054641844861809
$ws.Range("E15").Value = 14.14238458851022
$ws.Range("F15").Value = 24.67478063778865
$ws.Range("G15").Value = 29.58850993308851
$ws.Range("H15").Value = 14.17900979249226
$ws.Range("K15").Value = 16.16313886643621
$ws.Range("L15").Value = 9.220095708814128
$ws.Range("N15").Value = 16.74172518702457
$ws.Range("O15").Value = 21.7595624974332
$ws.Range("C16").Value = 13.86281064120784
$ws.Range("D16").Value = 5.024997555849699
$ws.Range("E16").Value = 14.08035590468997
$ws.Range("F16").Value = 24.6399226099384
$ws.Range("G16").Value = 29.52317300861554
$ws.Range("H16").Value = 14.20077787250577
$ws.Range("K16").Value = 15.83017314312406
$ws.Range("L16").Value = 9.214185386093863
$ws.Range("N16").Value = 16.75140791556827
$ws.Range("O16").Value = 21.77609738959851
$ws.Range("C17").Value = 13.80785150485068
$ws.Range("D17").Value = 5.006633799927702
$ws.Range("E17").Value = 14.0430058182466
$ws.Range("F17").Value = 24.6205320942804
$ws.Range("G17").Value = 29.4860289315516
$ws.Range("H17").Value = 14.21487643019959
$ws.Range("K17").Value = 15.6225097037802
$ws.Range("L17").Value = 9.210959322393933
$ws.Range("N17").Value = 16.75781296365631
$ws.Range("O17").Value = 21.78779279956597
$ws.Range("C18").Value = 13.77646207509111
$ws.Range("D18").Value = 4.996007516336476
$ws.Range("E18").Value = 14.02178679543178
$ws.Range("F18").Value = 24.61011386783585
$ws.Range("G18").Value = 29.46574838058003
$ws.Range("H18").Value = 14.22325881839068
$ws.Range("K18").Value = 15.50185157697933
$ws.Range("L18").Value = 9.209251100958548
$ws.Range("N18").Value = 16.76166813463998
$ws.Range("O18").Value = 21.79508898298745
$ws.Range("C19").Value = 13.76587313647999
$ws.Range("D19").Value = 4.992398797433858
$ws.Range("E19").Value = 14.0146482216762
$ws.Range("F19").Value = 24.6067127888923
$ws.Range("G19").Value = 29.45906821455927
$ws.Range("H19").Value = 14.22614385203823
$ws.Range("K19").Value = 15.46079338387016
$ws.Range("L19").Value = 9.208698068453636
$ws.Range("N19").Value = 16.76300285085679
$ws.Range("O19").Value = 21.79765702343697
$ws.Range("C20").Value = 13.81367928209636
$ws.Range("D20").Value = 5.008595303372132
$ws.Range("E20").Value = 14.04695462527127
$ws.Range("F20").Value = 24.62252025226534
$ws.Range("G20").Value = 29.48987090361952
$ws.Range("H20").Value = 14.21334732423031
$ws.Range("K20").Value = 15.64474243060885
$ws.Range("L20").Value = 9.211287504100349
$ws.Range("N20").Value = 16.75711342727397
$ws.Range("O20").Value = 21.78648886040182
$ws.Range("C21").Value = 13.97688051343533
$ws.Range("D21").Value = 5.062239226590828
$ws.Range("E21").Value = 14.15861172937292
$ws.Range("F21").Value = 24.68439490490362
$ws.Range("G21").Value = 29.6062833531402
$ws.Range("H21").Value = 14.17362016420108
$ws.Range("K21").Value = 16.24802679351917
$ws.Range("L21").Value = 9.221744854595961
$ws.Range("N21").Value = 16.73936883855565
$ws.Range("O21").Value = 21.75577402019567
$ws.Range("C22").Value = 14.08530670237806
$ws.Range("D22").Value = 5.096698575482566
$ws.Range("E22").Value = 14.23381439174631
$ws.Range("F22").Value = 24.73132989336792
$ws.Range("G22").Value = 29.69192312242381
$ws.Range("H22").Value = 14.15009687132411
$ws.Range("K22").Value = 16.63079624934205
$ws.Range("L22").Value = 9.229881824245597
$ws.Range("N22").Value = 16.72929012095226
$ws.Range("O22").Value = 21.74078938142679
$ws.Range("C23").Value = 14.02728630654506
$ws.Range("D23").Value = 5.07836334280819
$ws.Range("E23").Value = 14.19347938861925
$ws.Range("F23").Value = 24.70568792363751
$ws.Range("G23").Value = 29.64534616461804
$ws.Range("H23").Value = 14.16242802952333
$ws.Range("K23").Value = 16.42758907743086
$ws.Range("L23").Value = 9.225420279079959
$ws.Range("N23").Value = 16.73453066977022
$ws.Range("O23").Value = 21.74832042904218
$ws.Range("C24").Value = 13.81104389595593
$ws.Range("D24").Value = 5.007708721348222
$ws.Range("E24").Value = 14.04516857700104
$ws.Range("F24").Value = 24.62161913277142
$ws.Range("G24").Value = 29.48813060139945
$ws.Range("H24").Value = 14.21403777044873
$ws.Range("K24").Value = 15.63469496150606
$ws.Range("L24").Value = 9.211138676764307
$ws.Range("N24").Value = 16.7574291492651
$ws.Range("O24").Value = 21.78707658870774
$ws.Range("C25").Value = 13.58638142065535
$ws.Range("D25").Value = 4.92907078888104
$ws.Range("E25").Value = 13.89533872135613
$ws.Range("F25").Value = 24.55903688779938
$ws.Range("G25").Value = 29.36000488081408
$ws.Range("H25").Value = 14.28015015434787
$ws.Range("K25").Value = 14.73349325137896
$ws.Range("L25").Value = 9.201370121354675
$ws.Range("N25").Value = 16.78867573990696
$ws.Range("O25").Value = 21.85063016943189

Write-Output "Updated loading_percent values for case Case_0_237 (380 kV)"
